# "Body size measurements and sample updates"
# Fill in newly-recorded sample counts on the tracker sheet and leave the
# selection where the user's cursor ended up after entering them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Site ARM - Vertical sample count
$ws.Range("B2").Value = 30

# Site EBM - Surface sample count
$ws.Range("C3").Value = 3

# Site GMS - Vertical sample count
$ws.Range("B4").Value = 30

# Site PIM - Vertical sample count
$ws.Range("B5").Value = 15

# Cursor ends on B6 after data entry
$ws.Range("B6").Select()
